$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = 68.17709991409968
$ws.Range("E1").Value = 72.14331646379448
$ws.Range("F1").Value = 73.05243123138962
$ws.Range("G1").Value = 42.88452623217619

$ws.Range("D2").Value = 80.53147070105067
$ws.Range("E2").Value = 85.07462686567165
$ws.Range("F2").Value = 85.60849030440302
$ws.Range("G2").Value = 54.64641705294696

$ws.Range("D3").Value = 74.45049886913807
$ws.Range("E3").Value = 79.75264751006495
$ws.Range("F3").Value = 80.44781521322891
$ws.Range("G3").Value = 51.88594645843765

$ws.Range("D4").Value = 59.52251958002392
$ws.Range("E4").Value = 71.64330274225603
$ws.Range("F4").Value = 67.56447684451997
$ws.Range("G4").Value = 45.05974583798825

$ws.Range("D5").Value = 69.38626482441448
$ws.Range("E5").Value = 76.66666666666667
$ws.Range("F5").Value = 83.40054994615512
$ws.Range("G5").Value = 47.64263970045284

$ws.Range("D6").Value = 71.4383148407932
$ws.Range("E6").Value = 55.00000000000001
$ws.Range("F6").Value = 61.14347932559139
$ws.Range("G6").Value = 23.55577734366527

$ws.Range("D7").Value = 72.66193058683237
$ws.Range("E7").Value = 68.42105263157895
$ws.Range("F7").Value = 75.06667551845975
$ws.Range("G7").Value = 38.00636203072422

$ws.Range("D8").Value = 74.61180520482702
$ws.Range("E8").Value = 70.58823529411765
$ws.Range("F8").Value = 68.25713015681517
$ws.Range("G8").Value = 40.72825938824338

$ws.Range("D9").Value = 42.81399470571772
$ws.Range("E9").Value = 70
$ws.Range("F9").Value = 62.93083254194365
$ws.Range("G9").Value = 41.55106204495092
